# osi-and-tcp.pptx - "Add files via upload"
#
# The OSI/TCP-IP comparison table on the slide gets its "Layer Number"
# column updated:
#   - The "7" / "6" / "5" cells (Application / Presentation / Session
#     rows) are merged into a single vertically-spanning cell and
#     relabeled "L7(Firewall)".
#   - The remaining single-row layer numbers are relabeled to include
#     the matching network-device hint: "4" -> "L4(NAT)",
#     "3" -> "L3(Router)", "2" -> "L2(Switch)", "1" -> "L1".

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

# "Layer Number" is the 3rd column of the table.
$layerCol = 3

$cell7 = $tbl.Cell(2, $layerCol)   # "7" - Application row
$cell6 = $tbl.Cell(3, $layerCol)   # "6" - Presentation row
$cell5 = $tbl.Cell(4, $layerCol)   # "5" - Session row

# Merge the three cells vertically (mirrors the existing HTTP / SSH,FTP...
# columns that already span these same three rows) and relabel the
# resulting cell.
$cell7.Merge($cell5)
$cell7.Shape.TextFrame.TextRange.Text = "L7(Firewall)"

# The now-covered cells keep no visible text.
$cell6.Shape.TextFrame.TextRange.Text = ""
$cell5.Shape.TextFrame.TextRange.Text = ""

# Relabel the remaining (still single-row) layer numbers.
$tbl.Cell(5, $layerCol).Shape.TextFrame.TextRange.Text = "L4(NAT)"
$tbl.Cell(6, $layerCol).Shape.TextFrame.TextRange.Text = "L3(Router)"
$tbl.Cell(7, $layerCol).Shape.TextFrame.TextRange.Text = "L2(Switch)"
$tbl.Cell(8, $layerCol).Shape.TextFrame.TextRange.Text = "L1"
